# create_forecast_ad_hoc/inputs_outputs.xlsx
# "DP: create_forecast_ad_hoc - make promoteres dynamic"
#
# Repoints the hard-coded network (W:\...) paths used by the ad-hoc
# forecast script to the author's local working-copy paths, switches the
# scenario name from the "with project" run to the "without project" run,
# and bumps v_date to the matching version.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# B2: location of the base-forecast tool
$ws.Range("B2").Value = "C:\Users\dpere\Documents\JTMT\forecast_git\create_forecast_basic\current"

# B3: client/project location
$ws.Range("B3").Value = "C:\Users\dpere\Documents\JTMT\Projects\תחזיות_דמוגרפיות\קבצי עבודה\142_מתחם_אנגל\בהת"

# B4: forecast_version scenario name
$ws.Range("B4").Value = "without_project"

# B5: v_date
$ws.Range("B5").Value = 240129

# B6: base-forecast output location by version
$ws.Range("B6").Value = "C:\Users\dpere\Documents\JTMT\forecast_by_version\V4\BASE_YEAR"

# Row 3 picked up a slightly taller height in the saved file.
$ws.Rows("3:3").RowHeight = 14.5

# Active cell ends on the v_date input this time.
[void]$ws.Range("B5").Select()
